$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4083.2307
$ws.Range("I15").Value = 4083.2307
$ws.Range("K15").Value = 12249.6921
$ws.Range("M15").Value = -12080.6921
$ws.Range("H97").Value = 1953.72
$ws.Range("J97").Value = 1953.72
$ws.Range("L97").Value = 5861.16
$ws.Range("N97").Value = -6853.16
$ws.Range("H125").Value = 4942.5557
$ws.Range("I125").Value = 898
$ws.Range("K125").Value = 8082
$ws.Range("M125").Value = -5622
$ws.Range("H137").Value = 4765.1377
$ws.Range("J137").Value = 2316.4443
$ws.Range("L137").Value = 6949.3329
$ws.Range("N137").Value = -12049.3329
$ws.Range("H138").Value = 2530.776
$ws.Range("J138").Value = 2803.1353
$ws.Range("L138").Value = 8409.4059
$ws.Range("N138").Value = -18689.4059

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1504888.8
$ws.Range("I2").Value = 1681616.6
$ws.Range("J2").Value = 2702
$ws.Range("K2").Value = 1681616.6
$ws.Range("L2").Value = 2702
$ws.Range("M2").Value = -1681503.6
$ws.Range("N2").Value = -2928
$ws.Range("H32").Value = 5994.25
$ws.Range("I32").Value = 5478.96
$ws.Range("J32").Value = 10288.333
$ws.Range("K32").Value = 5478.96
$ws.Range("L32").Value = 10288.333
$ws.Range("M32").Value = -5191.96
$ws.Range("N32").Value = -10862.333
$ws.Range("H37").Value = 9475.25
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H61").Value = 4600.2915
$ws.Range("I61").Value = 5764.4546
$ws.Range("J61").Value = 3615.2307
$ws.Range("K61").Value = 5764.4546
$ws.Range("L61").Value = 3615.2307
$ws.Range("M61").Value = -5552.4546
$ws.Range("N61").Value = -4039.2307
$ws.Range("H74").Value = 2340.0588
$ws.Range("I74").Value = 1921.6786
$ws.Range("K74").Value = 1921.6786
$ws.Range("M74").Value = -1047.6786
$ws.Range("H77").Value = 2340.0588
$ws.Range("I77").Value = 1921.6786
$ws.Range("K77").Value = 9608.393
$ws.Range("M77").Value = -5240.393
$ws.Range("H102").Value = 3268.15
$ws.Range("I102").Value = 2365.8462
$ws.Range("J102").Value = 4943.857
$ws.Range("K102").Value = 2365.8462
$ws.Range("L102").Value = 4943.857
$ws.Range("M102").Value = -743.8462
$ws.Range("N102").Value = -8187.857
$ws.Range("H110").Value = 2977.1
$ws.Range("I110").Value = 3276.3333
$ws.Range("K110").Value = 3276.3333
$ws.Range("M110").Value = -1231.3333
$ws.Range("H116").Value = 1504888.8
$ws.Range("I116").Value = 1681616.6
$ws.Range("J116").Value = 2702
$ws.Range("K116").Value = 1681616.6
$ws.Range("L116").Value = 2702
$ws.Range("M116").Value = -1679322.6
$ws.Range("N116").Value = -7290
$ws.Range("H132").Value = 3701.7334
$ws.Range("I132").Value = 3271.3845
$ws.Range("K132").Value = 9814.1535
$ws.Range("M132").Value = -7284.1535
$ws.Range("H136").Value = 4600.2915
$ws.Range("I136").Value = 5764.4546
$ws.Range("J136").Value = 3615.2307
$ws.Range("K136").Value = 17293.3638
$ws.Range("L136").Value = 10845.6921
$ws.Range("M136").Value = -14743.3638
$ws.Range("N136").Value = -15945.6921

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1504888.8
$ws.Range("I3").Value = 1681616.6
$ws.Range("J3").Value = 2702
$ws.Range("K3").Value = 1681616.6
$ws.Range("L3").Value = 2702
$ws.Range("M3").Value = -1681502.6
$ws.Range("N3").Value = -2930
$ws.Range("H102").Value = 42557.07
$ws.Range("I102").Value = 65159.8
$ws.Range("K102").Value = 65159.8
$ws.Range("M102").Value = -61914.8
$ws.Range("H134").Value = 1829
$ws.Range("I134").Value = 1414.8148
$ws.Range("K134").Value = 4244.4444
$ws.Range("M134").Value = -1709.4444

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3489.5557
$ws.Range("I31").Value = 1943.3334
$ws.Range("K31").Value = 1943.3334
$ws.Range("M31").Value = -1648.3334
$ws.Range("H34").Value = 3489.5557
$ws.Range("I34").Value = 1943.3334
$ws.Range("K34").Value = 1943.3334
$ws.Range("M34").Value = -1741.3334
$ws.Range("H107").Value = 1109.95
$ws.Range("I107").Value = 956.1818
$ws.Range("J107").Value = 1297.8889
$ws.Range("K107").Value = 956.1818
$ws.Range("L107").Value = 1297.8889
$ws.Range("M107").Value = 963.8182
$ws.Range("N107").Value = -5137.8889
$ws.Range("H120").Value = 29899
$ws.Range("J120").Value = 29899
$ws.Range("L120").Value = 29899
$ws.Range("N120").Value = -37157
$ws.Range("H121").Value = 64108.668
$ws.Range("J121").Value = 64108.668
$ws.Range("L121").Value = 64108.668
$ws.Range("N121").Value = -66728.66800000001
$ws.Range("H134").Value = 5294
$ws.Range("I134").Value = 5242.875
$ws.Range("K134").Value = 15728.625
$ws.Range("M134").Value = -13193.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 283.0909
$ws.Range("I14").Value = 283.0909
$ws.Range("K14").Value = 849.2727
$ws.Range("M14").Value = -676.2727
$ws.Range("H56").Value = 14879.682
$ws.Range("I56").Value = 14879.682
$ws.Range("K56").Value = 14879.682
$ws.Range("M56").Value = -14349.682
$ws.Range("H113").Value = 571.5714
$ws.Range("I113").Value = 685.25
$ws.Range("J113").Value = 420
$ws.Range("K113").Value = 2055.75
$ws.Range("L113").Value = 1260
$ws.Range("M113").Value = 114.25
$ws.Range("N113").Value = -5600
$ws.Range("H121").Value = 1265.7693
$ws.Range("J121").Value = 1619.6666
$ws.Range("L121").Value = 4858.9998
$ws.Range("N121").Value = -7478.9998
$ws.Range("H131").Value = 1840223.8
$ws.Range("I131").Value = 3268782.5
$ws.Range("J131").Value = 3505.5715
$ws.Range("K131").Value = 9806347.5
$ws.Range("L131").Value = 10516.7145
$ws.Range("M131").Value = -9801307.5
$ws.Range("N131").Value = -20596.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2790.2307
$ws.Range("I7").Value = 2652.5557
$ws.Range("J7").Value = 3100
$ws.Range("K7").Value = 2652.5557
$ws.Range("L7").Value = 3100
$ws.Range("M7").Value = -2540.5557
$ws.Range("N7").Value = -3324
$ws.Range("H22").Value = 862.9048
$ws.Range("I22").Value = 834
$ws.Range("K22").Value = 834
$ws.Range("M22").Value = -539
$ws.Range("H27").Value = 862.9048
$ws.Range("I27").Value = 834
$ws.Range("K27").Value = 834
$ws.Range("M27").Value = -727
$ws.Range("H122").Value = 5657.6
$ws.Range("I122").Value = 7337.857
$ws.Range("J122").Value = 4187.375
$ws.Range("K122").Value = 22013.571
$ws.Range("L122").Value = 12562.125
$ws.Range("M122").Value = -19563.571
$ws.Range("N122").Value = -17462.125
$ws.Range("H126").Value = 2790.2307
$ws.Range("I126").Value = 2652.5557
$ws.Range("J126").Value = 3100
$ws.Range("K126").Value = 7957.6671
$ws.Range("L126").Value = 9300
$ws.Range("M126").Value = -5487.6671
$ws.Range("N126").Value = -14240
$ws.Range("H136").Value = 12935.857
$ws.Range("J136").Value = 5555
$ws.Range("L136").Value = 16665
$ws.Range("N136").Value = -21765

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6256.6
$ws.Range("I122").Value = 4257.7896
$ws.Range("K122").Value = 12773.3688
$ws.Range("M122").Value = -10323.3688
$ws.Range("H126").Value = 1819.3846
$ws.Range("I126").Value = 1787
$ws.Range("J126").Value = 1997.5
$ws.Range("K126").Value = 5361
$ws.Range("L126").Value = 5992.5
$ws.Range("M126").Value = -2891
$ws.Range("N126").Value = -10932.5
$ws.Range("H132").Value = 3234.5574
$ws.Range("I132").Value = 3181.5085
$ws.Range("J132").Value = 4799.5
$ws.Range("K132").Value = 9544.5255
$ws.Range("L132").Value = 14398.5
$ws.Range("M132").Value = -7014.5255
$ws.Range("N132").Value = -19458.5
